# Add the "InvalidLogin" test-data sheet (negative login test case) after
# the existing "ValidLogin" sheet, matching the author's commit:
#   "Created TestData for invalid login TC
#    Developed Script for invalid login TC"

$wb = $excel.ActiveWorkbook

# Locate the existing sheet and insert the new one right after it so tab
# order becomes ValidLogin, InvalidLogin.
$ws1 = $wb.Worksheets.Item("ValidLogin")
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "InvalidLogin"

# Same header row as ValidLogin, but with an invalid (bad) credential pair.
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "abcd"
$ws2.Range("B2").Value = "xyz"

# Match the zoom level / selection of the recorded sheet view for the new
# (now active) sheet.
$ws2.Application.ActiveWindow.Zoom = 220
$ws2.Range("E5").Select() | Out-Null
